$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / title rich text updates ---
# A8 shared string: "Volume 30   Number  4" -> "...Number  5" (only last run's digit changes)
$a8 = $ws.Range("A8")
$a8.Characters(21, 1).Text = "5"

# C9 shared string: "Report Covering the Week  1/23/2023  Through  1/29/2023"
# -> "...1/30/2023  Through  2/5/2023" (2nd run "1/23/2023"->"1/30/2023", 4th run "1/29/2023"->"2/5/2023")
$c9 = $ws.Range("C9")
$c9.Characters(27, 9).Text = "1/30/2023"
$c9.Characters(47, 9).Text = "2/5/2023"


# --- Crime statistics table updates (rows 14-30) ---
# --- Cells whose type/style changes (requires format copy + value set) ---
$ws.Range("N14").Copy()
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range("M14").Value = -100

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F23").PasteSpecial(-4122)

$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4122)

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)

# --- Cells with simple value updates (style/class unchanged) ---
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = -40
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 0
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 20
$ws.Range("F16").Value = 27
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = 12.5
$ws.Range("I16").Value = 41
$ws.Range("J16").Value = 35
$ws.Range("K16").Value = 17.142857142857
$ws.Range("L16").Value = 64
$ws.Range("M16").Value = -16.326530612244
$ws.Range("N16").Value = -75
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 71.428571428571
$ws.Range("F17").Value = 49
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = 48.484848484848
$ws.Range("I17").Value = 54
$ws.Range("J17").Value = 42
$ws.Range("K17").Value = 28.571428571428
$ws.Range("L17").Value = 10.204081632653
$ws.Range("M17").Value = 22.727272727272
$ws.Range("N17").Value = -20.588235294117
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 14
$ws.Range("H18").Value = 7.692307692307
$ws.Range("I18").Value = 23
$ws.Range("J18").Value = 18
$ws.Range("K18").Value = 27.777777777777
$ws.Range("L18").Value = 27.777777777777
$ws.Range("M18").Value = -50
$ws.Range("N18").Value = -89.449541284403
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 8.333333333333
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = -23.214285714285
$ws.Range("I19").Value = 52
$ws.Range("J19").Value = 70
$ws.Range("K19").Value = -25.714285714285
$ws.Range("L19").Value = -11.864406779661
$ws.Range("M19").Value = 1.960784313725
$ws.Range("N19").Value = -46.938775510204
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 62.5
$ws.Range("I20").Value = 38
$ws.Range("J20").Value = 27
$ws.Range("K20").Value = 40.74074074074
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = 111.111111111111
$ws.Range("N20").Value = -75.796178343949
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 27.586206896551
$ws.Range("F21").Value = 161
$ws.Range("G21").Value = 146
$ws.Range("H21").Value = 10.273972602739
$ws.Range("I21").Value = 211
$ws.Range("J21").Value = 197
$ws.Range("K21").Value = 7.106598984771
$ws.Range("L21").Value = 21.965317919075
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = -70.406732117812
$ws.Range("G22").Value = 2
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 81
$ws.Range("E24").Value = -59.259259259259
$ws.Range("F24").Value = 162
$ws.Range("G24").Value = 254
$ws.Range("H24").Value = -36.220472440944
$ws.Range("I24").Value = 207
$ws.Range("J24").Value = 305
$ws.Range("K24").Value = -32.131147540983
$ws.Range("L24").Value = 84.821428571428
$ws.Range("M24").Value = 51.094890510948
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 112.5
$ws.Range("F25").Value = 55
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = 7.843137254901
$ws.Range("I25").Value = 62
$ws.Range("J25").Value = 66
$ws.Range("K25").Value = -6.060606060606
$ws.Range("L25").Value = 21.56862745098
$ws.Range("M25").Value = -6.060606060606
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 9
$ws.Range("H26").Value = -66.666666666666
$ws.Range("I26").Value = 4
$ws.Range("J26").Value = 10
$ws.Range("K26").Value = -60
$ws.Range("L26").Value = 33.333333333333
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 350
$ws.Range("I27").Value = 12
$ws.Range("K27").Value = 200
$ws.Range("L27").Value = 50
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = -66.666666666666
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = -71.428571428571
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = -66.666666666666
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -66.666666666666
$ws.Range("G30").Value = 1